$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlVAlignCenter = -4108
$xlHAlignCenter = -4108

# ------------------------------------------------------------------
# 1. New row: "aug_1.0" run (trained with the fully augmented dataset)
# ------------------------------------------------------------------
$ws.Range("A4").Value = "Baseline3DConvNet"
$ws.Range("B4").Value = "aug_1.0"
$ws.Range("C4").Value = 1

# ------------------------------------------------------------------
# 2. New "Notes" column (F): widen it, then add header + notes text
# ------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 40.666666666666664

$ws.Range("F1").Value = "Notes"

$notesText = "1) Added dropout in the fully connected layer" + [char]10 + "2) Reduced learning rate from 1e-4 to 5e-5" + [char]10 + "3) train loss and accuracy are estimated on the original train dataset, not on th augment dataset"
$ws.Range("F4").Value = $notesText

# ------------------------------------------------------------------
# 3. Finish the "aug_0.5" row (row 3) with its measured results
#    and a note about how the numbers were computed
# ------------------------------------------------------------------
$ws.Range("D3").Value = 0.846797
$ws.Range("E3").Value = 0.858333
$ws.Range("F3").Value = "Train loss and accuracy are estimated at the end of epoch on augumented dataset"

# ------------------------------------------------------------------
# 4. Keep the accuracy number format consistent across every row,
#    including a few more blank rows reserved under the table
#    (D/E columns only)
# ------------------------------------------------------------------
$ws.Range("D2:E7").NumberFormat = "0.0000"

# ------------------------------------------------------------------
# 5. Row heights: header + the two wrapped note rows
# ------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 16
$ws.Rows.Item(3).RowHeight = 32
$ws.Rows.Item(4).RowHeight = 64

# ------------------------------------------------------------------
# 6. Formatting pass: vertically center everything, keep the existing
#    horizontal centering on the numeric columns, and wrap the new
#    Notes column text
# ------------------------------------------------------------------

# Header row (bold)
$ws.Range("A1:B1").VerticalAlignment = $xlVAlignCenter
$ws.Range("C1:E1").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C1:E1").VerticalAlignment = $xlVAlignCenter
$ws.Range("F1").VerticalAlignment = $xlVAlignCenter
$ws.Range("F1").WrapText = $true

# Model / variant columns (left aligned text, vertical center only)
$ws.Range("A2:B4").VerticalAlignment = $xlVAlignCenter

# frac_sample_aug column (center)
$ws.Range("C2:C4").HorizontalAlignment = $xlHAlignCenter
$ws.Range("C2:C4").VerticalAlignment = $xlVAlignCenter

# Train/valid accuracy columns (including the trailing blank rows)
$ws.Range("D2:E7").HorizontalAlignment = $xlHAlignCenter
$ws.Range("D2:E7").VerticalAlignment = $xlVAlignCenter

# Notes column data cells: vertical center + wrap
$ws.Range("F3:F4").VerticalAlignment = $xlVAlignCenter
$ws.Range("F3:F4").WrapText = $true

# ------------------------------------------------------------------
# 7. Leave the selection where the author last left it
# ------------------------------------------------------------------
[void]$ws.Range("D7").Select()

Write-Host "done"
